$wb = $excel.ActiveWorkbook

# --- Ticket 50: add the "varStatus" demo sheet, placed after "Immaterial" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$newSheet.Name = "varStatus"

# Header row (X, Y, startX, endX, stepX, startY, endY, stepY)
$newSheet.Range("A1").Value = "X"
$newSheet.Range("B1").Value = "Y"
$newSheet.Range("D1").Value = "endX"
$newSheet.Range("C1").Value = "startX"
$newSheet.Range("E1").Value = "stepX"
$newSheet.Range("F1").Value = "startY"
$newSheet.Range("G1").Value = "endY"
$newSheet.Range("H1").Value = "stepY"

# Template row with the nested jt:for / varStatus expressions
$newSheet.Range("B2").Value = "`${y}"
$newSheet.Range("A2").Value = '<jt:for var="x" start="1" end="5" varStatus="xs"><jt:for var="y" start="5" end="${x}" step="-1" varStatus="ys">${x}'
$newSheet.Range("C2").Value = '${xs.start}'
$newSheet.Range("F2").Value = '${ys.start}'
$newSheet.Range("G2").Value = '${ys.end}'
$newSheet.Range("H2").Value = '${ys.step}</jt:for></jt:for>'
$newSheet.Range("D2").Value = '${xs.end}'
$newSheet.Range("E2").Value = '${xs.step}'

# Row 1 header style: bold font + blue fill + thin border, reusing the same
# look already used for headers on the other sheets (copy format so the
# existing style slot is reused instead of minting a near-duplicate one).
$srcSheet = $wb.Worksheets.Item("Multiplication")
$srcSheet.Range("A1").Copy()
$newSheet.Range("A1:H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 2 style: plain font, no fill, thin border (new style).
$newSheet.Range("A2:H2").Borders.LineStyle = 1

# Match portrait page setup used by the other sheets.
$newSheet.PageSetup.Orientation = 1
